# Applies the "Improved documentation in Test Suite" commit to the
# TestSuiteSubmission workbook.
#
# Summary of the real (content-level) changes:
#  1. "Constant Samples"!C9   : "alt cell correct" -> "alt cell(given in key) correct"
#  2. "Constant Samples"!C11  : (new) "Negative but correct, hence zero"
#  3. "Constant Samples"!C12  : (new) "Negative (Penalty)"
#  4. "Formula Samples"!C4    : old "ROUNDUP ... wrong" note -> new "ROUNDUP ... correct" note
#  5. "Formula Samples"!C11   : "pre-req right" -> "pre-req wrong"
#  6. Per-sheet active cell / selection changes on several worksheets.
#  7. Active-sheet changes to "SoftFormula Samples" (tabSelected moves there).

$wb = $excel.ActiveWorkbook

# --- 1) Constant Samples: new / updated notes -----------------------------
$wsConstant = $wb.Worksheets.Item("Constant Samples")
$wsConstant.Range("C9").Value = "alt cell(given in key) correct"
$wsConstant.Range("C11").Value = "Negative but correct, hence zero"
$wsConstant.Range("C12").Value = "Negative (Penalty)"

# --- 2) Formula Samples: updated notes -------------------------------------
$wsFormula = $wb.Worksheets.Item("Formula Samples")
$wsFormula.Range("C4").Value = "> This cell will be considered correct since even though ROUNDUP is a custom formula, and it multiplies against 0.01 instead of doing what the key does: dividing against 100. But with SymPy similiarity check, it succeeds"
$wsFormula.Range("C11").Value = "pre-req wrong"

# --- 3) Selection / active-cell updates on each sheet -----------------------
$wsCheck = $wb.Worksheets.Item("Check Samples")
$wsCheck.Range("A4").Select()

$wsConstant.Range("C5").Select()

$wsFormula.Range("C24").Select()

$wsRelative = $wb.Worksheets.Item("Relative Samples")
$wsRelative.Range("C16").Select()

$wsMinimum = $wb.Worksheets.Item("Minimum Work")
$wsMinimum.Range("E15").Select()

# Select "SoftFormula Samples" last so it becomes the active / tabSelected sheet.
$wsSoft = $wb.Worksheets.Item("SoftFormula Samples")
$wsSoft.Range("B2").Select()
